$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the paragraph index (1-based) whose Range.Text equals $text,
# scanning forward from $startAt. Returns -1 if not found.
# ---------------------------------------------------------------------------
function Find-ParaIndex($doc, $text, $startAt) {
    $count = $doc.Paragraphs.Count
    for ($i = $startAt; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Helper: the source document has a handful of paragraphs whose short label
# ("Cart:", "Wishlist:", ...) was typed as two runs separated by a
# spell-check proofErr pair (<w:proofErr spellStart/>Word<w:proofErr
# spellEnd/>:). Re-typing the label collapses it back into a single,
# proofErr-free run while preserving the paragraph's list/number formatting
# and the run's underline formatting.
#
# We rebuild the paragraph from a *neighbour* paragraph that already has the
# exact paragraph formatting we want (same list level / style / bold state)
# so the freshly created paragraph mark inherits the correct wPr, then we
# delete the old (proofErr-laden) paragraph and fill in the new one's text.
#
# $anchorText   : exact Range.Text (including trailing CR) of the neighbour
#                 paragraph to clone formatting from.
# $anchorBefore : $true  -> new paragraph is inserted *after* the anchor
#                 $false -> new paragraph is inserted *before* the anchor
# ---------------------------------------------------------------------------
function Merge-LabelRuns($doc, $targetText, $newLabel, $anchorText, $anchorBefore, $searchStart) {
    $targetIdx = Find-ParaIndex $doc $targetText $searchStart
    if ($targetIdx -eq -1) {
        throw ("paragraph not found: " + $targetText)
    }

    if ($anchorBefore) {
        # anchor paragraph sits immediately before the target
        $anchorIdx = $targetIdx - 1
    } else {
        # anchor paragraph sits immediately after the target
        $anchorIdx = $targetIdx + 1
    }
    $anchor = $doc.Paragraphs.Item($anchorIdx)
    if ($anchor.Range.Text -ne $anchorText) {
        throw ("anchor paragraph mismatch: expected [" + $anchorText + "] got [" + $anchor.Range.Text + "]")
    }

    $target = $doc.Paragraphs.Item($targetIdx)
    $full = $doc.Range($target.Range.Start, $target.Range.End)
    $full.Delete()

    if ($anchorBefore) {
        # anchor index unchanged (it is before the deleted paragraph)
        $anchor = $doc.Paragraphs.Item($anchorIdx)
        $anchor.Range.InsertParagraphAfter()
        $createdIdx = $anchorIdx + 1
    } else {
        # anchor shifted up by one once the target paragraph was removed
        $anchor = $doc.Paragraphs.Item($anchorIdx - 1)
        $anchor.Range.InsertParagraphBefore()
        $createdIdx = $anchorIdx - 1
    }
    $created = $doc.Paragraphs.Item($createdIdx)
    $created.Range.Text = $newLabel

    # Hand back the index the rebuilt paragraph now occupies, so callers can
    # resume their next search just past it (the new paragraph's text would
    # otherwise re-match an identical $targetText on a later call).
    return $createdIdx
}

# ---------------------------------------------------------------------------
# Helper: delete a whole paragraph (including its paragraph mark) whose
# Range.Text matches $text exactly.
# ---------------------------------------------------------------------------
function Remove-Paragraph($doc, $text, $searchStart) {
    $idx = Find-ParaIndex $doc $text $searchStart
    if ($idx -eq -1) {
        throw ("paragraph not found: " + $text)
    }
    $p = $doc.Paragraphs.Item($idx)
    $full = $doc.Range($p.Range.Start, $p.Range.End)
    $full.Delete()
}

# ---------------------------------------------------------------------------
# 1) Collapse the "Cart"/"spellStart"/"spellEnd"/":" run-split back into a
#    single "Cart:" run, for every affected paragraph (in document order).
# ---------------------------------------------------------------------------
$pos = Merge-LabelRuns $d "Cart:`r" "Cart:" "Product:`r" $true 1
$next = $pos + 1
$pos = Merge-LabelRuns $d "Cart:`r" "Cart:" "Product Page:`r" $true $next
$next = $pos + 1
$pos = Merge-LabelRuns $d "Wishlist:`r" "Wishlist:" "All Modals:`r" $true $next
$next = $pos + 1
$pos = Merge-LabelRuns $d "Cart:`r" "Cart:" "Homepage:`r" $false $next

# ---------------------------------------------------------------------------
# 2) Remove the three stray QA-note paragraphs called out in the diff.
# ---------------------------------------------------------------------------
Remove-Paragraph $d "Drop down sizes are completely incorrect. Use lightbox example given by Ray. `r" 1
Remove-Paragraph $d "Add the VIEW ALL…. Button to each drop down. Will link to that category top level cat.`r" 1
Remove-Paragraph $d "Placing order text in the loading modal is black instead of white. Run by Michael on how we can change this text.`r" 1
